# Updates the "cryptos" price/volume sheet with refreshed values and, for a
# handful of rows, swapped coin rows (name/link/price/volume moved together).
# Column D ("Price") values are prefixed with a leading apostrophe so Excel
# keeps them as literal text (matching the source inlineStr cells) instead of
# re-parsing them as floating point numbers and losing their exact formatting
# (trailing zeros, multi-dot thousand separators, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''63.638.95'
$ws.Range('E2').Value = '  -1.12%  '
$ws.Range('D3').Value = '''2.628.36'
$ws.Range('E3').Value = '  +0.41%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '''575.84'
$ws.Range('E5').Value = '  +0.07%  '
$ws.Range('D6').Value = '''154.75'
$ws.Range('E6').Value = '  -0.78%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').Value = '''0.624'
$ws.Range('E8').Value = '  -2.72%  '
$ws.Range('D9').Value = '''2.625.03'
$ws.Range('E9').Value = '  +0.39%  '
$ws.Range('D10').Value = '''0.117'
$ws.Range('E10').Value = '  -3.23%  '
$ws.Range('D11').Value = '''5.80'
$ws.Range('E11').Value = '  +0.69%  '
$ws.Range('D12').Value = '''0.381'
$ws.Range('E12').Value = '  -1.79%  '
$ws.Range('E13').Value = '  +0.63%  '
$ws.Range('D14').Value = '''28.32'
$ws.Range('E14').Value = '  +0.32%  '
$ws.Range('D15').Value = '''3.096.88'
$ws.Range('E15').Value = '  +0.56%  '
$ws.Range('D16').Value = '''0.0000183'
$ws.Range('E16').Value = '  -0.86%  '
$ws.Range('D17').Value = '''63.532.63'
$ws.Range('E17').Value = '  -1.06%  '
$ws.Range('D18').Value = '''2.641.13'
$ws.Range('E18').Value = '  +1.43%  '
$ws.Range('D19').Value = '''12.10'
$ws.Range('E19').Value = '  -0.62%  '
$ws.Range('D20').Value = '''7.57'
$ws.Range('E20').Value = '  +3.60%  '
$ws.Range('D21').Value = '''4.52'
$ws.Range('E21').Value = '  -2.81%  '
$ws.Range('D22').Value = '''344.14'
$ws.Range('E22').Value = '  +0.39%  '
$ws.Range('E23').Value = '  +0.48%  '
$ws.Range('D24').Value = '''67.86'
$ws.Range('E24').Value = '  -0.17%  '
$ws.Range('D25').Value = '''1.87'
$ws.Range('E25').Value = '  +9.74%  '
$ws.Range('D26').Value = '''0.0000108'
$ws.Range('E26').Value = '  -1.68%  '
$ws.Range('D27').Value = '''597.04'
$ws.Range('E27').Value = '  +8.42%  '
$ws.Range('D28').Value = '''9.22'
$ws.Range('E28').Value = '  -0.88%  '
$ws.Range('D29').Value = '''1.60'
$ws.Range('E29').Value = '  +4.74%  '
$ws.Range('B30').Value = 'Aptos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D30').Value = '''7.93'
$ws.Range('E30').Value = '  +0.86%  '
$ws.Range('D31').Value = '''0.161'
$ws.Range('E31').Value = '  -0.46%  '
$ws.Range('B32').Value = 'Binance-PegBSC-USD'
$ws.Range('C32').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D32').Value = '''0.999'
$ws.Range('E32').Value = '  -0.19%  '
$ws.Range('D33').Value = '''2.05'
$ws.Range('E33').Value = '  -0.51%  '
$ws.Range('D34').Value = '''1.74'
$ws.Range('E34').Value = '  +1.17%  '
$ws.Range('D35').Value = '''6.60'
$ws.Range('E35').Value = '  +4.04%  '
$ws.Range('D36').Value = '''5.37'
$ws.Range('E36').Value = '  +1.57%  '
$ws.Range('D37').Value = '''0.401'
$ws.Range('E37').Value = '  -1.70%  '
$ws.Range('B38').Value = 'EthereumClassic'
$ws.Range('C38').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D38').Value = '''19.73'
$ws.Range('E38').Value = '  -0.95%  '
$ws.Range('B39').Value = 'FirstDigitalUSD'
$ws.Range('C39').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D39').Value = '''0.998'
$ws.Range('E39').Value = '  -0.06%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').Value = '''1.90'
$ws.Range('E40').Value = '  -0.29%  '
$ws.Range('B41').Value = 'Monero'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D41').Value = '''151.33'
$ws.Range('E41').Value = '  +0.13%  '
$ws.Range('E42').Value = '  -0.07%  '
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').Value = '''2.52'
$ws.Range('E43').Value = '  +5.34%  '
$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D44').Value = '''41.67'
$ws.Range('E44').Value = '  -0.24%  '
$ws.Range('D45').Value = '''158.41'
$ws.Range('E45').Value = '  +0.54%  '
$ws.Range('D46').Value = '''24.22'
$ws.Range('E46').Value = '  +7.58%  '
$ws.Range('D47').Value = '''3.90'
$ws.Range('E47').Value = '  -1.45%  '
$ws.Range('D48').Value = '''0.0587'
$ws.Range('E48').Value = '  -1.93%  '
$ws.Range('D49').Value = '''0.629'
$ws.Range('E49').Value = '  -0.30%  '
$ws.Range('D50').Value = '''0.0999'
$ws.Range('E50').Value = '  -0.77%  '
$ws.Range('D51').Value = '''0.0248'
$ws.Range('E51').Value = '  -0.27%  '
